# daily auto push: 2026-02-23 03:18 UTC
#
# Insert one new data row at row 839 ("2026/02/23" / 月 / 7 / 21), which
# pushes the existing rows 839:880 down to 840:881. The inserted row's
# A/B columns duplicate the date/weekday already present in row 838, so
# we copy that row down (this preserves the "General"-formatted text
# value for the date cell instead of letting Excel's value-setter
# re-interpret "2026/02/23" as a date serial) and then correct the two
# numeric columns for the newly inserted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 838 ("2026/02/23" / 月 / ...) as the new row 839, shifting
# the old rows 839-880 down to 840-881.
$ws.Rows.Item(838).Copy()
$ws.Rows.Item(839).Insert()

# Fix up the time/ranking values for the newly inserted row.
$ws.Range("C839").Value = 7
$ws.Range("D839").Value = 21
